# Cosmetic changes: insert two new blank columns immediately before the
# existing "C" column, so the original column C (and everything to its
# right) shifts two columns right to become column E -- carrying its
# per-cell styles and its custom column width along with it. Populate
# the two freshly inserted columns (new C and new D) with the new
# header text / placeholder values, and restore the header text that
# used to live in column B into the new column D (its new home),
# giving column B its own new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at C:D. Former column C (and D, E, ...)
# shift two columns to the right; former column C becomes column E,
# preserving its contents, per-cell styles and column width.
$ws.Columns("C:D").Insert()

# Row 1 header labels.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Fill the new C and D columns (rows 2-27) with the same placeholder
# value ("UN") used throughout the rest of the table.
$ws.Range("C2:D27").Value = "UN"

# Give the two new columns (and re-assert it on the column that used
# to be C, now E) the same custom width used by the rest of the
# table's data columns.
$ws.Columns("C:E").ColumnWidth = 7.1666666666667
